$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2025-05-28 07:20:15", "uuuuuuuuuuu", "Import - Nouveau produit", 25, 0, 25, "2481023879"),
    @("2025-05-28 07:21:35", "uuuuuuuuuuu", "Sortie", 8, 25, 17, "2481023879"),
    @("2025-05-28 07:21:51", "uuuuuuuuuuu", "Entrée", 25, 17, 42, "2481023879"),
    @("2025-05-28 07:22:02", "uuuuuuuuuuu", "Sortie", 8, 42, 34, "2481023879"),
    @("2025-05-28 07:22:07", "uuuuuuuuuuu", "Sortie", 5, 34, 29, "2481023879"),
    @("2025-05-28 07:22:16", "uuuuuuuuuuu", "Entrée", 23, 29, 52, "2481023879"),
    @("2025-05-28 07:22:24", "uuuuuuuuuuu", "Sortie", 12, 52, 40, "2481023879")
)

$rowIndex = 20
foreach ($rowData in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rowData[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rowData[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rowData[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rowData[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rowData[4]
    $ws.Cells.Item($rowIndex, 6).Value = $rowData[5]
    $refCell = $ws.Cells.Item($rowIndex, 7)
    $refCell.NumberFormat = "@"
    $refCell.Value = $rowData[6]
    $refCell.Style = "Normal"
    $rowIndex++
}
